$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add "producto comodin" - wildcard product used so stock is not updated for it
$ws.Range("B36").Value = "Agustina"
$ws.Range("C36").Value = 1
$ws.Range("C36").NumberFormat = "0%"

# Update the active selection shown when the file was last saved
$ws.Range("B41").Select()
